$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "2022" column (S) mirroring the existing "2021" column (R):
# copy each R-column cell's formatting down into the matching S-column cell,
# then fill in the new 2022 figures.
# ---------------------------------------------------------------------------

# Bring each row's format across from column R to column S (rows 4-43, the
# data block; header rows 1-3 get no new cell, matching the source change).
$ws.Range("R4:R43").Copy()
$ws.Range("S4:S43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New 2022 values for the data rows (rows 8 and 36 are section headers with
# no numeric figure, so they stay blank).
$ws.Cells.Item(4, 19).Value = 2022
$ws.Cells.Item(5, 19).Value = 4.9000000000000004
$ws.Cells.Item(6, 19).Value = 6.1
$ws.Cells.Item(7, 19).Value = 4
$ws.Cells.Item(9, 19).Value = 6.1
$ws.Cells.Item(10, 19).Value = 12.4
$ws.Cells.Item(11, 19).Value = 3.2
$ws.Cells.Item(12, 19).Value = 10.8
$ws.Cells.Item(13, 19).Value = 14.6
$ws.Cells.Item(14, 19).Value = 8.5
$ws.Cells.Item(15, 19).Value = 5.5
$ws.Cells.Item(16, 19).Value = 7.1
$ws.Cells.Item(17, 19).Value = 4.4000000000000004
$ws.Cells.Item(18, 19).Value = 5.8
$ws.Cells.Item(19, 19).Value = 11.6
$ws.Cells.Item(20, 19).Value = 3.1
$ws.Cells.Item(21, 19).Value = 1.5
$ws.Cells.Item(22, 19).Value = 2.2999999999999998
$ws.Cells.Item(23, 19).Value = 1
$ws.Cells.Item(24, 19).Value = 2.2999999999999998
$ws.Cells.Item(25, 19).Value = 3.3
$ws.Cells.Item(26, 19).Value = 1.6
$ws.Cells.Item(27, 19).Value = 4.5999999999999996
$ws.Cells.Item(28, 19).Value = 4.4000000000000004
$ws.Cells.Item(29, 19).Value = 4.7
$ws.Cells.Item(30, 19).Value = 4
$ws.Cells.Item(31, 19).Value = 3.2
$ws.Cells.Item(32, 19).Value = 4.7
$ws.Cells.Item(33, 19).Value = 2.6
$ws.Cells.Item(34, 19).Value = 3.3
$ws.Cells.Item(35, 19).Value = 2.2000000000000002
$ws.Cells.Item(37, 19).Value = 13.2
$ws.Cells.Item(38, 19).Value = 7.5
$ws.Cells.Item(39, 19).Value = 4.0999999999999996
$ws.Cells.Item(40, 19).Value = 4.3
$ws.Cells.Item(41, 19).Value = 2.6
$ws.Cells.Item(42, 19).Value = 1

# Row 43 footnote cell: same "…" marker as column R.
$ws.Cells.Item(43, 19).Value = $ws.Cells.Item(43, 18).Value()

# The two section-header rows (8 and 36) get a distinct bold-italic style
# rather than a plain copy of column R's formatting.
$ws.Range("G1").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("S8").Font.Bold = $true
$ws.Range("S8").Font.Italic = $true
$ws.Range("S8").Value = $null

$ws.Range("G1").Copy()
$ws.Range("S36").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("S36").Font.Bold = $true
$ws.Range("S36").Font.Italic = $true
$ws.Range("S36").Value = $null

# Move the active selection the way the author left it.
$ws.Range("T12").Select()
